$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "HK_R_acc_G"
$ws.Range("A2").Value = 71.663619744058508
$ws.Range("A3").Value = 71.343692870201096
$ws.Range("A4").Value = 71.709323583180989
$ws.Range("A5").Value = 70.292504570383912
$ws.Range("A6").Value = 69.149908592321751
$ws.Range("A7").Value = 70.338208409506393
$ws.Range("A8").Value = 68.510054844606955
$ws.Range("A9").Value = 69.652650822669102
$ws.Range("A10").Value = 69.058500914076788
$ws.Range("A11").Value = 69.012797074954292
$ws.Range("A12").Value = 70.292504570383912
$ws.Range("A13").Value = 70.566727605118828
$ws.Range("A14").Value = 69.744058500914079
$ws.Range("A15").Value = 72.486288848263257
$ws.Range("A16").Value = 69.9725776965265
$ws.Range("A17").Value = 69.698354661791598
$ws.Range("A18").Value = 70.795246800731263
$ws.Range("A19").Value = 70.658135283363805
$ws.Range("A20").Value = 75.091407678244977
$ws.Range("A21").Value = 72.760511882998173
$ws.Range("A22").Value = 75
$ws.Range("A23").Value = 68.967093235831811
$ws.Range("A24").Value = 67.093235831809878
$ws.Range("A25").Value = 69.378427787934186
$ws.Range("A26").Value = 69.515539305301644
$ws.Range("A27").Value = 70.566727605118828
$ws.Range("A28").Value = 69.78976234003656
$ws.Range("A29").Value = 73.308957952468006
$ws.Range("A30").Value = 73.263254113345525
$ws.Range("A31").Value = 72.806215722120655
$ws.Range("A32").Value = 68.235831809872025
$ws.Range("A33").Value = 69.378427787934186
$ws.Range("A34").Value = 71.06946983546618
$ws.Range("A35").Value = 70.292504570383912
$ws.Range("A36").Value = 70.292504570383912
$ws.Range("A37").Value = 72.577696526508234
$ws.Range("A38").Value = 69.515539305301644
$ws.Range("A39").Value = 69.78976234003656
$ws.Range("A40").Value = 70.42961608775137
$ws.Range("A41").Value = 69.78976234003656
$ws.Range("A42").Value = 70.201096892138935
$ws.Range("A43").Value = 69.881170018281537
$ws.Range("A44").Value = 69.469835466179163
$ws.Range("A45").Value = 70.201096892138935
$ws.Range("A46").Value = 71.800731261425966
$ws.Range("A47").Value = 69.012797074954292
$ws.Range("A48").Value = 68.921389396709316
$ws.Range("A49").Value = 70.978062157221217
